# feat : real email service
#
# Appointments sheet: widen the patient_phone column, promote the
# member_id / group_number values on the existing last row (11) from
# text to numbers, and append a brand-new appointment as row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- widen column P (patient_phone, column 16) from 19 to 28 characters ---
# ColumnWidth is specified in "characters of the Normal font" and Excel
# stores the OOXML <col width> value with a constant +5/6 padding baked
# in, so request 5/6 less than the target to land exactly on 28.
$ws.Columns.Item(16).ColumnWidth = 28 - (5/6)

# --- row 11: member_id / group_number were stored as text "12345";
#     they should be real numbers now ---
$ws.Cells.Item(11, 11).Value = 12345
$ws.Cells.Item(11, 12).Value = 12345

# --- row 12: new appointment record ---
# K12/L12 (member_id/group_number) and R12/S12 (date_of_birth /
# appointment_date) look numeric/date-like but must stay plain text,
# matching every other text cell in the sheet, so force a text number
# format before writing them and restore General afterwards.
$ws.Range("K12:L12").NumberFormat = "@"
$ws.Range("R12:S12").NumberFormat = "@"

$ws.Cells.Item(12, 1).Value = "APT_20250906_034542_001"
$ws.Cells.Item(12, 2).Value = "PAT_056"
$ws.Cells.Item(12, 3).Value = "Shreyansh Bhatter"
$ws.Cells.Item(12, 4).Value = "Dr. Aish"
$ws.Cells.Item(12, 5).Value = "Banjara Hills"
# F12/G12 (date, time) are blank for this record
$ws.Cells.Item(12, 8).Value = 30
$ws.Cells.Item(12, 9).Value = "confirmed"
$ws.Cells.Item(12, 10).Value = "Cigna"
$ws.Cells.Item(12, 11).Value = "12345"
$ws.Cells.Item(12, 12).Value = "12345"
$ws.Cells.Item(12, 13).Value = "2025-09-06T03:45:42.554058"
$ws.Cells.Item(12, 14).Value = $false
$ws.Cells.Item(12, 15).Value = 0
$ws.Cells.Item(12, 16).Value = "shreyanshs070700@gmail.com"
$ws.Cells.Item(12, 17).Value = "(701) 368-4370"
$ws.Cells.Item(12, 18).Value = "07/07/2004"
$ws.Cells.Item(12, 19).Value = "2025-09-08"
$ws.Cells.Item(12, 20).Value = "13:00"
$ws.Cells.Item(12, 21).Value = "2025-09-06T03:45:42.554058"
$ws.Cells.Item(12, 22).Value = $true
$ws.Cells.Item(12, 23).Value = "2025-09-06T03:45:47.094612"

# restore General number format on the cells we temporarily forced to text
$ws.Range("K12:L12").NumberFormat = "general"
$ws.Range("R12:S12").NumberFormat = "general"

# match the rest of the table's row styling (thin border, left/center
# aligned) for the whole new row, including the two blank cells
$ws.Range("A12:W12").Borders.LineStyle = 1
$ws.Range("A12:W12").HorizontalAlignment = -4131
$ws.Range("A12:W12").VerticalAlignment = -4108
